# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (G) values for each game row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 0
    11 = 0
    12 = 1
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 2
    18 = 0
    19 = 2
    20 = 1
    21 = 0
    22 = 1
    23 = 1
    24 = 0
    25 = 0
    26 = 2
    27 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
